# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.470.19"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "4.022.90"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +18.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.765"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000321"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.37%  "
$ws.Range("D14").Value = "4.661.71"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "4.018.18"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.132"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "72.135.96"
$ws.Range("E20").Value = "  +3.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "102.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "668.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "0.0₃0859"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.428"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.150"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0483"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.159"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.00%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.71%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.02%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000266"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.96%  "
